# Apply attendance updates to the sheet.
# Mapping of row -> cell updates (column letter = new value)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @{ G = 1; H = 1 }
    4  = @{ D = 1; E = 1 }
    5  = @{ D = 2; E = 1; F = 1 }
    6  = @{ D = 1; E = 1 }
    7  = @{ H = 1 }
    8  = @{ H = 1 }
    9  = @{ H = 1 }
    10 = @{ D = 1; E = 1 }
    11 = @{ H = 1 }
    12 = @{ H = 1 }
    13 = @{ D = 1; E = 1 }
    14 = @{ H = 1 }
    15 = @{ H = 1 }
    16 = @{ H = 1 }
    17 = @{ D = 1; E = 1 }
    18 = @{ H = 1 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
